$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.734.16"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").Value = "1.803.28"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.20"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.05"
$ws.Range("E8").Value = "  +7.14%  "
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0673"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "2.063.15"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.36"
$ws.Range("E13").Value = "  +15.05%  "
$ws.Range("D14").Value = "1.807.33"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "34.749.31"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.30"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.70"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "256.90"
$ws.Range("E19").Value = "  +3.13%  "
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.52"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.26"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.90"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.56"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.17"
$ws.Range("E27").Value = "  +4.57%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("E34").Value = "  +11.23%  "
$ws.Range("D35").Value = "1.461.36"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0191"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.638"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.75"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.907"
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("E44").Value = "  +6.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0509"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "1.962.44"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.09"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.72"
$ws.Range("E50").Value = "  +5.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.85"
$ws.Range("E51").Value = "  -0.96%  "
